$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H81").Value = 38328
$ws.Range("J81").Value = 38328
$ws.Range("L81").Value = 38328
$ws.Range("N81").Value = -40324

$ws.Range("H84").Value = 38328
$ws.Range("J84").Value = 38328
$ws.Range("L84").Value = 114984
$ws.Range("N84").Value = -124968

$ws.Range("H101").Value = 628.3333
$ws.Range("I101").Value = 217.14285
$ws.Range("J101").Value = 2067.5
$ws.Range("K101").Value = 651.4285500000001
$ws.Range("L101").Value = 6202.5
$ws.Range("M101").Value = 970.5714499999999
$ws.Range("N101").Value = -9446.5

$ws.Range("H137").Value = 493405.03
$ws.Range("I137").Value = 3161.7812
$ws.Range("J137").Value = 929176.8
$ws.Range("K137").Value = 9485.3436
$ws.Range("L137").Value = 2787530.4
$ws.Range("M137").Value = -6935.3436
$ws.Range("N137").Value = -2792630.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H95").Value = 48993.332
$ws.Range("J95").Value = 48993.332
$ws.Range("L95").Value = 48993.332
$ws.Range("N95").Value = -54485.332

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2527907
$ws.Range("I58").Value = 4134773.5
$ws.Range("J58").Value = 2831
$ws.Range("K58").Value = 4134773.5
$ws.Range("L58").Value = 2831
$ws.Range("M58").Value = -4134570.5
$ws.Range("N58").Value = -3237

$ws.Range("H136").Value = 2527907
$ws.Range("I136").Value = 4134773.5
$ws.Range("J136").Value = 2831
$ws.Range("K136").Value = 12404320.5
$ws.Range("L136").Value = 8493
$ws.Range("M136").Value = -12401770.5
$ws.Range("N136").Value = -13593

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 4112.697
$ws.Range("I3").Value = 2517.7222
$ws.Range("J3").Value = 6026.6665
$ws.Range("K3").Value = 7553.1666
$ws.Range("L3").Value = 18079.9995
$ws.Range("M3").Value = -7441.1666
$ws.Range("N3").Value = -18303.9995

$ws.Range("H56").Value = 171633.5
$ws.Range("I56").Value = 171633.5
$ws.Range("K56").Value = 171633.5
$ws.Range("M56").Value = -171103.5

$ws.Range("H133").Value = 3360
$ws.Range("I133").Value = 3360
$ws.Range("J133").Value = 0
$ws.Range("K133").Value = 10080
$ws.Range("L133").Value = 0
$ws.Range("M133").Value = -5020
$ws.Range("N133").ClearContents()

$ws.Range("H140").Value = 2374.8108
$ws.Range("I140").Value = 1835.9259
$ws.Range("J140").Value = 3829.8
$ws.Range("K140").Value = 5507.7777
$ws.Range("L140").Value = 11489.4
$ws.Range("M140").Value = -327.7776999999996
$ws.Range("N140").Value = -21849.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H19").Value = 19000
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 19000
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 19000
$ws.Range("M19").ClearContents()
$ws.Range("N19").Value = -19340

$ws.Range("H22").Value = 2099.6667
$ws.Range("I22").Value = 2099.6667
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 2099.6667
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -1804.6667
$ws.Range("N22").ClearContents()

$ws.Range("H27").Value = 2099.6667
$ws.Range("I27").Value = 2099.6667
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 2099.6667
$ws.Range("L27").Value = 0
$ws.Range("M27").Value = -1992.6667
$ws.Range("N27").ClearContents()

$ws.Range("H46").Value = 1320.2
$ws.Range("I46").Value = 1467.3334
$ws.Range("J46").Value = 1099.5
$ws.Range("K46").Value = 1467.3334
$ws.Range("L46").Value = 1099.5
$ws.Range("M46").Value = -1279.3334
$ws.Range("N46").Value = -1475.5

$ws.Range("H55").Value = 128.35294
$ws.Range("I55").Value = 128.36363
$ws.Range("K55").Value = 128.36363
$ws.Range("M55").Value = 44.63637

$ws.Range("H61").Value = 67375.5
$ws.Range("J61").Value = 35000
$ws.Range("L61").Value = 35000
$ws.Range("N61").Value = -35404

$ws.Range("H113").Value = 67375.5
$ws.Range("J113").Value = 35000
$ws.Range("L113").Value = 35000
$ws.Range("N113").Value = -39340

$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H117:N117").ClearContents()
$ws.Range("H118:N118").ClearContents()
$ws.Range("H119:N119").ClearContents()
$ws.Range("H120:N120").ClearContents()
$ws.Range("H122:N122").ClearContents()
$ws.Range("H123:N123").ClearContents()
$ws.Range("H124:N124").ClearContents()
$ws.Range("H125:N125").ClearContents()
$ws.Range("H126:N126").ClearContents()
$ws.Range("H127:N127").ClearContents()
$ws.Range("H128:N128").ClearContents()
$ws.Range("H129:N129").ClearContents()
$ws.Range("H130:N130").ClearContents()
$ws.Range("H131:N131").ClearContents()
$ws.Range("H132:N132").ClearContents()
$ws.Range("H133:N133").ClearContents()
$ws.Range("H134:N134").ClearContents()
$ws.Range("H135:N135").ClearContents()
$ws.Range("H137:N137").ClearContents()
$ws.Range("H138:N138").ClearContents()
$ws.Range("H139:N139").ClearContents()
$ws.Range("H140:N140").ClearContents()
$ws.Range("H141:N141").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H119:N119").ClearContents()
$ws.Range("H120:N120").ClearContents()
$ws.Range("H121:N121").ClearContents()
$ws.Range("H122:N122").ClearContents()
$ws.Range("H123:N123").ClearContents()
$ws.Range("H124:N124").ClearContents()
$ws.Range("H125:N125").ClearContents()
$ws.Range("H126:N126").ClearContents()
$ws.Range("H127:N127").ClearContents()
$ws.Range("H128:N128").ClearContents()
$ws.Range("H129:N129").ClearContents()
$ws.Range("H130:N130").ClearContents()
$ws.Range("H131:N131").ClearContents()
$ws.Range("H132:N132").ClearContents()
$ws.Range("H133:N133").ClearContents()
$ws.Range("H135:N135").ClearContents()
$ws.Range("H136:N136").ClearContents()
$ws.Range("H137:N137").ClearContents()
$ws.Range("H138:N138").ClearContents()
$ws.Range("H139:N139").ClearContents()
$ws.Range("H140:N140").ClearContents()
$ws.Range("H141:N141").ClearContents()
